$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A93").Value = "gen:10075"
$ws.Range("B93").Value = "contributor type"

$ws.Range("A94").Value = "gen:10076"
$ws.Range("B94").Value = "data collector"
$ws.Range("G94").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H94").Value = "gen:10075"

$ws.Range("A95").Value = "gen:10077"
$ws.Range("B95").Value = "data manager"
$ws.Range("G95").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H95").Value = "gen:10075"

$ws.Range("A96").Value = "gen:10078"
$ws.Range("B96").Value = "data curator"
$ws.Range("G96").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H96").Value = "gen:10075"

$ws.Range("A97").Value = "gen:10079"
$ws.Range("B97").Value = "sponsor"
$ws.Range("G97").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H97").Value = "gen:10075"

$ws.Range("A98").Value = "gen:10080"
$ws.Range("B98").Value = "distributor"
$ws.Range("G98").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H98").Value = "gen:10075"

$ws.Range("A99").Value = "gen:10081"
$ws.Range("B99").Value = "project leader"
$ws.Range("G99").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H99").Value = "gen:10075"

$ws.Range("A100").Value = "gen:10082"
$ws.Range("B100").Value = "project manager"
$ws.Range("G100").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H100").Value = "gen:10075"

$ws.Range("A101").Value = "gen:10083"
$ws.Range("B101").Value = "project member"
$ws.Range("G101").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H101").Value = "gen:10075"

$ws.Range("A102").Value = "gen:10084"
$ws.Range("B102").Value = "reseacher"
$ws.Range("G102").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H102").Value = "gen:10075"

$ws.Range("A103").Value = "gen:10085"
$ws.Range("B103").Value = "editor"
$ws.Range("G103").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H103").Value = "gen:10075"

$ws.Range("A104").Value = "gen:10086"
$ws.Range("B104").Value = "rights holder"
$ws.Range("G104").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H104").Value = "gen:10075"

$ws.Range("A105").Value = "gen:10087"
$ws.Range("B105").Value = "supervisor"
$ws.Range("G105").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H105").Value = "gen:10075"

$ws.Range("A106").Value = "gen:10088"
$ws.Range("B106").Value = "work package leader"
$ws.Range("G106").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H106").Value = "gen:10075"

$ws.Range("A107").Value = "gen:10089"
$ws.Range("B107").Value = "contact person"
$ws.Range("G107").Value = "https://schema.datacite.org/meta/kernel-4.4/doc/DataCite-MetadataKernel_v4.4.pdf"
$ws.Range("H107").Value = "gen:10075"
